$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:K1").NumberFormat = "@"

$ws.Cells.Item(1,1).Value = "Holton"
$ws.Cells.Item(1,2).Value = "Johnny"
$ws.Cells.Item(1,3).Value = "WR"
$ws.Cells.Item(1,4).Value = "2018-11-25"
$ws.Cells.Item(1,5).Value = "11"
$ws.Cells.Item(1,6).Value = "27.095"
$ws.Cells.Item(1,7).Value = "OAK"
$ws.Cells.Item(1,8).Value = "@"
$ws.Cells.Item(1,9).Value = "BAL"
$ws.Cells.Item(1,10).Value = "L 17-34"
$ws.Cells.Item(1,11).Value = ""
$ws.Cells.Item(1,12).Value = 0
